# Add 2022-Q4 data
# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet named "2022-Q4" right before the existing
#    "2022-Q3" sheet (so it becomes the 2nd tab, right after "总计").
# 2. Populate it with the per-fund holdings table for 2022-Q4.
# 3. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 right
#    after the header row, shifting the existing rows down, and renumber the
#    running index in column A.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create & position the new "2022-Q4" worksheet.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row (same headers used by every per-quarter sheet).
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Make column B (fund code, has leading zeros) and D..G text-typed
# (these number-formatted cells hold the values as text in the source
# data, e.g. "010610"/"47.74" rather than numeric 10610/47.74).
$q4.Range("B2:B11").NumberFormat = "@"
$q4.Range("D2:G11").NumberFormat = "@"

# Column A (running index) uses the same bold / bordered / centered style
# as the header row and as column A on every other per-quarter sheet.
$aCol = $q4.Range("A2:A11")
$aCol.Font.Bold = $true
$aCol.Borders.LineStyle = 1
$aCol.HorizontalAlignment = -4108
$aCol.VerticalAlignment = -4160

$q4Data = @(
    @(0, "010610", "上投摩根远见两年持有期混合",             "47.74", "91.22", "3.07", "1.4656", 10),
    @(1, "375010", "上投摩根中国优势混合A",                  "19.80", "87.18", "2.75", "0.5445", 8),
    @(2, "014261", "上投摩根沃享远见一年持有期混合A",         "12.10", "87.47", "2.98", "0.3606", 8),
    @(3, "630011", "华商主题精选混合",                        "3.43", "92.61", "4.00", "0.1372", 7),
    @(4, "013142", "华商乐享互联灵活配置混合C",                "5.21", "88.52", "2.30", "0.1198", 9),
    @(5, "001959", "华商乐享互联灵活配置混合A",                "4.50", "88.52", "2.30", "0.1035", 9),
    @(6, "015709", "上投摩根中国优势混合C",                    "2.00", "87.18", "2.75", "0.0550", 8),
    @(7, "233001", "大摩基础行业混合",                          "0.70", "78.60", "5.82", "0.0407", 6),
    @(8, "014262", "上投摩根沃享远见一年持有期混合C",           "0.67", "87.47", "2.98", "0.0200", 8),
    @(9, "510560", "国寿安保中证500ETF",                        "1.81", "99.21", "0.31", "0.0056", 5)
)

$rowNum = 2
foreach ($rec in $q4Data) {
    $q4.Cells.Item($rowNum, 1).Value = $rec[0]
    $q4.Cells.Item($rowNum, 2).Value = $rec[1]
    $q4.Cells.Item($rowNum, 3).Value = $rec[2]
    $q4.Cells.Item($rowNum, 4).Value = $rec[3]
    $q4.Cells.Item($rowNum, 5).Value = $rec[4]
    $q4.Cells.Item($rowNum, 6).Value = $rec[5]
    $q4.Cells.Item($rowNum, 7).Value = $rec[6]
    $q4.Cells.Item($rowNum, 8).Value = $rec[7]
    $rowNum = $rowNum + 1
}

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new 2022-Q4 row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Existing rows 2..8 hold: 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q2,
# 2021-Q1, 2020-Q4 (in that order). Read them out before overwriting so we
# can re-write them one row further down.
$existingRows = @()
for ($r = 2; $r -le 8; $r++) {
    $existingRows += ,@(
        $total.Cells.Item($r, 2).Value(),
        $total.Cells.Item($r, 3).Value(),
        $total.Cells.Item($r, 4).Value()
    )
}

# Give the new last row (row 9) the same formatting as the existing data
# rows before writing into it (column A needs the bold "index" style).
$total.Cells.Item(2, 1).Copy($total.Cells.Item(9, 1))

# Re-write rows 2..8 shifted down to rows 3..9, renumbering column A.
for ($i = 0; $i -lt $existingRows.Length; $i++) {
    $destRow = $i + 3
    $total.Cells.Item($destRow, 1).Value = $i + 1
    $total.Cells.Item($destRow, 2).Value = $existingRows[$i][0]
    $total.Cells.Item($destRow, 3).Value = $existingRows[$i][1]
    $total.Cells.Item($destRow, 4).Value = $existingRows[$i][2]
}

# Write the brand-new 2022-Q4 row into row 2.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 10
$total.Cells.Item(2, 4).Value = 2.85

Write-Output "Done: 2022-Q4 sheet added and summary sheet updated"
